# Apply cryptos list price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.380.42'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").Value = '3.527.37'
$ws.Range("E3").Value = '  -3.72%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'609.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -5.16%  '
$ws.Range("D6").Value = "'153.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("D7").Value = '3.521.72'
$ws.Range("E7").Value = '  -3.67%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'0.486"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("D11").Value = "'6.89"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("D12").Value = "'0.429"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("E13").Value = '  -4.04%  '
$ws.Range("D14").Value = '4.126.76'
$ws.Range("E14").Value = '  -3.64%  '
$ws.Range("D15").Value = "'31.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '3.528.24'
$ws.Range("E16").Value = '  -3.24%  '
$ws.Range("D17").Value = '67.434.85'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("D19").Value = "'6.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("D20").Value = "'15.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.25%  '
$ws.Range("D21").Value = "'452.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.07%  '
$ws.Range("D22").Value = "'9.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.49%  '
$ws.Range("D23").Value = "'0.641"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").Value = "'78.67"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("D25").Value = '3.674.57'
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").Value = "'10.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("D29").Value = "'8.32"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.80%  '
$ws.Range("D30").Value = "'2.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("D31").Value = "'1.67"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = "'25.94"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("D34").Value = "'1.90"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.43%  '
$ws.Range("D35").Value = "'0.158"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("D36").Value = "'6.20"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.78%  '
$ws.Range("D37").Value = '3.530.68'
$ws.Range("E37").Value = '  -3.39%  '
$ws.Range("D38").Value = "'7.96"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.27%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = "'176.41"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").Value = "'5.61"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("D43").Value = "'0.0878"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("D45").Value = "'0.892"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.49%  '
$ws.Range("D46").Value = "'29.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +9.13%  '
$ws.Range("D47").Value = "'45.72"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("D48").Value = "'2.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("E49").Value = '  -1.92%  '
$ws.Range("D50").Value = "'7.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("E51").Value = '  -3.50%  '
